$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 3
$ws.Range("H3").Value = 42880.832
$ws.Range("J3").Value = 42880.832
$ws.Range("L3").Value = 42880.832
$ws.Range("N3").Value = -43108.832
# Row 32
$ws.Range("H32").Value = 1422.4286
$ws.Range("I32").Value = 3000
$ws.Range("J32").Value = 1301.0769
$ws.Range("K32").Value = 3000
$ws.Range("L32").Value = 1301.0769
$ws.Range("M32").Value = -2674
$ws.Range("N32").Value = -1953.0769
# Row 51
$ws.Range("H51").Value = 2807.6924
$ws.Range("I51").Value = 1525
$ws.Range("J51").Value = 3377.7778
$ws.Range("K51").Value = 1525
$ws.Range("L51").Value = 3377.7778
$ws.Range("M51").Value = -1041
$ws.Range("N51").Value = -4345.7778
# Row 86
$ws.Range("H86").Value = 29415396
$ws.Range("I86").Value = 55558770
$ws.Range("J86").Value = 4100
$ws.Range("K86").Value = 55558770
$ws.Range("L86").Value = 4100
$ws.Range("M86").Value = -55557647
$ws.Range("N86").Value = -6346
# Row 89
$ws.Range("H89").Value = 29415396
$ws.Range("I89").Value = 55558770
$ws.Range("J89").Value = 4100
$ws.Range("K89").Value = 277793850
$ws.Range("L89").Value = 20500
$ws.Range("M89").Value = -277788234
$ws.Range("N89").Value = -31732
# Row 102
$ws.Range("H102").Value = 42880.832
$ws.Range("J102").Value = 42880.832
$ws.Range("L102").Value = 42880.832
$ws.Range("N102").Value = -49370.832
# Row 129
$ws.Range("H129").Value = 2068.6775
$ws.Range("J129").Value = 2995
$ws.Range("L129").Value = 8985
$ws.Range("N129").Value = -18985
# Row 131
$ws.Range("H131").Value = 5471.1816
$ws.Range("I131").Value = 851.7778
$ws.Range("J131").Value = 8669.23
$ws.Range("K131").Value = 2555.3334
$ws.Range("L131").Value = 26007.69
$ws.Range("M131").Value = 2484.6666
$ws.Range("N131").Value = -36087.69
# Row 132
$ws.Range("H132").Value = 7356524.5
$ws.Range("I132").Value = 3261
$ws.Range("J132").Value = 62506000
$ws.Range("K132").Value = 9783
$ws.Range("L132").Value = 187518000
$ws.Range("M132").Value = -7253
$ws.Range("N132").Value = -187523060
# Row 135
$ws.Range("H135").Value = 1091
$ws.Range("I135").Value = 1118.0322
$ws.Range("J135").Value = 971.2857
$ws.Range("K135").Value = 10062.2898
$ws.Range("L135").Value = 8741.5713
$ws.Range("M135").Value = -7527.2898
$ws.Range("N135").Value = -13811.5713
# Row 136
$ws.Range("H136").Value = 105999.5
$ws.Range("J136").Value = 105999.5
$ws.Range("L136").Value = 105999.5
$ws.Range("N136").Value = -116199.5
# Row 137
$ws.Range("H137").Value = 8696653
$ws.Range("I137").Value = 661.0833
$ws.Range("J137").Value = 18183190
$ws.Range("K137").Value = 1983.2499
$ws.Range("L137").Value = 54549570
$ws.Range("M137").Value = 566.7501
$ws.Range("N137").Value = -54554670
# Row 138
$ws.Range("H138").Value = 8548953
$ws.Range("I138").Value = 11112858
$ws.Range("J138").Value = 2603.3333
$ws.Range("K138").Value = 33338574
$ws.Range("L138").Value = 7809.999899999999
$ws.Range("M138").Value = -33333434
$ws.Range("N138").Value = -18089.9999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 7693488.5
$ws.Range("I61").Value = 9616614
$ws.Range("J61").Value = 986.8461
$ws.Range("K61").Value = 9616614
$ws.Range("L61").Value = 986.8461
$ws.Range("M61").Value = -9616402
$ws.Range("N61").Value = -1410.8461
# Row 74
$ws.Range("H74").Value = 13891878
$ws.Range("I74").Value = 25002768
$ws.Range("J74").Value = 3263.8125
$ws.Range("K74").Value = 25002768
$ws.Range("L74").Value = 3263.8125
$ws.Range("M74").Value = -25001894
$ws.Range("N74").Value = -5011.8125
# Row 77
$ws.Range("H77").Value = 13891878
$ws.Range("I77").Value = 25002768
$ws.Range("J77").Value = 3263.8125
$ws.Range("K77").Value = 125013840
$ws.Range("L77").Value = 16319.0625
$ws.Range("M77").Value = -125009472
$ws.Range("N77").Value = -25055.0625
# Row 132
$ws.Range("H132").Value = 6099490
$ws.Range("I132").Value = 7814318
$ws.Range("J132").Value = 2323.5557
$ws.Range("K132").Value = 23442954
$ws.Range("L132").Value = 6970.6671
$ws.Range("M132").Value = -23440424
$ws.Range("N132").Value = -12030.6671
# Row 136
$ws.Range("H136").Value = 7693488.5
$ws.Range("I136").Value = 9616614
$ws.Range("J136").Value = 986.8461
$ws.Range("K136").Value = 28849842
$ws.Range("L136").Value = 2960.5383
$ws.Range("M136").Value = -28847292
$ws.Range("N136").Value = -8060.5383

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 59
$ws.Range("H59").Value = 90000
$ws.Range("J59").Value = 90000
$ws.Range("L59").Value = 90000
$ws.Range("N59").Value = -91694
# Row 134
$ws.Range("H134").Value = 2230.1516
$ws.Range("I134").Value = 1338.7446
$ws.Range("J134").Value = 4435.2104
$ws.Range("K134").Value = 4016.2338
$ws.Range("L134").Value = 13305.6312
$ws.Range("M134").Value = -1481.2338
$ws.Range("N134").Value = -18375.6312

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 9014665
$ws.Range("I31").Value = 6527.76
$ws.Range("J31").Value = 27781618
$ws.Range("K31").Value = 6527.76
$ws.Range("L31").Value = 27781618
$ws.Range("M31").Value = -6232.76
$ws.Range("N31").Value = -27782208
# Row 34
$ws.Range("H34").Value = 9014665
$ws.Range("I34").Value = 6527.76
$ws.Range("J34").Value = 27781618
$ws.Range("K34").Value = 6527.76
$ws.Range("L34").Value = 27781618
$ws.Range("M34").Value = -6325.76
$ws.Range("N34").Value = -27782022
# Row 86
$ws.Range("H86").Value = 3050.5293
$ws.Range("I86").Value = 2736.3635
$ws.Range("J86").Value = 3626.5
$ws.Range("K86").Value = 2736.3635
$ws.Range("L86").Value = 3626.5
$ws.Range("M86").Value = -1613.3635
$ws.Range("N86").Value = -5872.5
# Row 89
$ws.Range("H89").Value = 3050.5293
$ws.Range("I89").Value = 2736.3635
$ws.Range("J89").Value = 3626.5
$ws.Range("K89").Value = 13681.8175
$ws.Range("L89").Value = 18132.5
$ws.Range("M89").Value = -8065.817499999999
$ws.Range("N89").Value = -29364.5
# Row 132
$ws.Range("H132").Value = 17859762
$ws.Range("I132").Value = 20835568
$ws.Range("J132").Value = 4928.5
$ws.Range("K132").Value = 62506704
$ws.Range("L132").Value = 14785.5
$ws.Range("M132").Value = -62504174
$ws.Range("N132").Value = -19845.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 33
$ws.Range("H33").Value = 25000068
$ws.Range("I33").Value = 28571512
$ws.Range("J33").Value = 20000044
$ws.Range("K33").Value = 171429072
$ws.Range("L33").Value = 120000264
$ws.Range("M33").Value = -171428789
$ws.Range("N33").Value = -120000830
# Row 134
$ws.Range("H134").Value = 3444.6572
$ws.Range("I134").Value = 1878.6364
$ws.Range("J134").Value = 6094.846
$ws.Range("K134").Value = 5635.9092
$ws.Range("L134").Value = 18284.538
$ws.Range("M134").Value = -565.9092000000001
$ws.Range("N134").Value = -28424.538

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 3261.3547
$ws.Range("I102").Value = 3347.0715
$ws.Range("J102").Value = 2461.3333
$ws.Range("K102").Value = 3347.0715
$ws.Range("L102").Value = 2461.3333
$ws.Range("M102").Value = -1725.0715
$ws.Range("N102").Value = -5705.3333

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 139
$ws.Range("H139").Value = 56716.668
$ws.Range("J139").Value = 56716.668
$ws.Range("L139").Value = 56716.668
$ws.Range("N139").Value = -66996.66800000001
